$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 550
$ws.Range("B4").Value = 180
$ws.Range("B5").Value = 23
$ws.Range("B6").Value = 150
$ws.Range("B8").Value = 400
